$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.396.78'
$ws.Range('E2').Value = '  -1.13%  '
$ws.Range('D3').Value = '1.565.85'
$ws.Range('E3').Value = '  -1.19%  '
$ws.Range('E4').Value = '  -0.05%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '208.27'
$c.Style = "Normal"
$ws.Range('E5').Value = '  +0.47%  '
$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '0.499'
$c.Style = "Normal"
$ws.Range('E6').Value = '  -0.91%  '
$ws.Range('E7').Value = '  -0.01%  '
$c = $ws.Range('D8')
$c.NumberFormat = "@"
$c.Value = '21.89'
$c.Style = "Normal"
$ws.Range('E8').Value = '  -2.07%  '
$c = $ws.Range('D9')
$c.NumberFormat = "@"
$c.Value = '0.249'
$c.Style = "Normal"
$ws.Range('E9').Value = '  -2.12%  '
$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '0.0591'
$c.Style = "Normal"
$ws.Range('E10').Value = '  -0.04%  '
$c = $ws.Range('D11')
$c.NumberFormat = "@"
$c.Value = '0.0868'
$c.Style = "Normal"
$ws.Range('E11').Value = '  -0.04%  '
$ws.Range('D12').Value = '1.786.65'
$ws.Range('E12').Value = '  -1.23%  '
$ws.Range('D13').Value = '1.567.79'
$ws.Range('E13').Value = '  -0.82%  '
$ws.Range('E14').Value = '  -1.23%  '
$ws.Range('E15').Value = '  -2.94%  '
$c = $ws.Range('D16')
$c.NumberFormat = "@"
$c.Value = '63.33'
$c.Style = "Normal"
$ws.Range('D17').Value = '27.383.87'
$ws.Range('E17').Value = '  -1.06%  '
$ws.Range('B18').Value = 'BitcoinCash'
$ws.Range('C18').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$c = $ws.Range('D18')
$c.NumberFormat = "@"
$c.Value = '212.22'
$c.Style = "Normal"
$ws.Range('E18').Value = '  -2.56%  '
$ws.Range('B19').Value = 'ShibaInu'
$ws.Range('C19').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D19').Value = '0.0₃0689'
$ws.Range('E19').Value = '  -0.73%  '
$ws.Range('E21').Value = '  -0.03%  '
$ws.Range('E22').Value = '  -1.28%  '
$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '9.53'
$c.Style = "Normal"
$ws.Range('E23').Value = '  -0.42%  '
$ws.Range('E24').Value = '  +1.28%  '
$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '153.38'
$c.Style = "Normal"
$ws.Range('E25').Value = '  -0.18%  '
$ws.Range('E26').Value = '  -0.06%  '
$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '6.73'
$c.Style = "Normal"
$ws.Range('E27').Value = '  +0.35%  '
$ws.Range('E28').Value = '  -0.63%  '
$ws.Range('E29').Value = '  -2.01%  '
$ws.Range('E30').Value = '  -0.26%  '
$ws.Range('E31').Value = '  +0.97%  '
$ws.Range('D33').Value = '1.362.97'
$ws.Range('E33').Value = '  -1.06%  '
$ws.Range('E34').Value = '  -0.17%  '
$ws.Range('E35').Value = '  +1.26%  '
$c = $ws.Range('D36')
$c.NumberFormat = "@"
$c.Value = '0.972'
$c.Style = "Normal"
$ws.Range('E36').Value = '  +0.14%  '
$ws.Range('E37').Value = '  -0.26%  '
$ws.Range('E38').Value = '  +0.72%  '
$c = $ws.Range('D39')
$c.NumberFormat = "@"
$c.Value = '0.532'
$c.Style = "Normal"
$ws.Range('E39').Value = '  -1.09%  '
$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '0.822'
$c.Style = "Normal"
$ws.Range('E40').Value = '  +0.67%  '
$ws.Range('E41').Value = '  -0.02%  '
$c = $ws.Range('D42')
$c.NumberFormat = "@"
$c.Value = '0.972'
$c.Style = "Normal"
$ws.Range('E42').Value = '  -0.36%  '
$ws.Range('E43').Value = '  +0.13%  '
$c = $ws.Range('D44')
$c.NumberFormat = "@"
$c.Value = '64.04'
$c.Style = "Normal"
$ws.Range('E44').Value = '  +0.52%  '
$ws.Range('E45').Value = '  +1.03%  '
$ws.Range('E46').Value = '  -1.17%  '
$ws.Range('D47').Value = '1.700.91'
$ws.Range('E47').Value = '  -1.05%  '
$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '85.56'
$c.Style = "Normal"
$ws.Range('E48').Value = '  -2.49%  '
$ws.Range('D49').Value = '0.0₇0994'
$ws.Range('E49').Value = '  -0.43%  '
$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '0.0956'
$c.Style = "Normal"
$ws.Range('E50').Value = '  -1.89%  '
$ws.Range('E51').Value = '  -0.65%  '
